$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (styles + row height) of row 22 onto the two new
# rows (23, 24) before filling in values, so the new cells inherit the
# same per-column styles used throughout the table (style 1 on col A,
# style 2 on cols C/D, default style elsewhere).
$ws.Range("A22:G22").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122)
$ws.Range("A24:G24").PasteSpecial(-4122)

$ws.Rows.Item(23).RowHeight = 13.8
$ws.Rows.Item(24).RowHeight = 13.8

# Row 23 - new review
$ws.Range("A23").Value = "com.hamxa.shaynachim"
$ws.Range("B23").Value = "bitcoin"
$ws.Range("C23").Value = "levenglickdov@gmail.com"
$ws.Range("D23").Value = "rotemzinger3@gmail.com"
$ws.Range("E23").Value = "27/5/2019 15:59"
$ws.Range("F23").Value = "put some salt and eat it like a pro"
$ws.Range("G23").Value = "yes"

# Row 24 - new review
$ws.Range("A24").Value = "com.hamxa.shaynachim"
$ws.Range("B24").Value = "bitcoin"
$ws.Range("C24").Value = "kassachale437@gmail.com"
$ws.Range("D24").Value = "levenglickdov@gmail.com"
$ws.Range("E24").Value = "27/5/2019 15:59"
$ws.Range("F24").Value = "delicious bitcoin app"
$ws.Range("G24").Value = "yes"

# Recovery column (D23) links to its matching recovery-email review, the
# same way existing rows such as D21 do.
$ws.Hyperlinks.Add($ws.Range("D23"), "mailto:rotemzinger3@gmail.com", "", "", "rotemzinger3@gmail.com")

# Update selection / scroll to match the saved view state.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 2
$ws.Range("G25").Select()
